$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A144").Value = "Плюс дополнительно покупал по мелочи в ближайших магазинах, что обязательно нужно:"

$ws.Range("A145").Value = "-"
$ws.Range("B145").Value = "медикаменты, лекарства"

$ws.Range("A146").Value = "-"
$ws.Range("B146").Value = "тапочки (полиуретан)"

$ws.Range("A147").Value = "-"
$ws.Range("B147").Value = "швейные принадлежности (нитки, иголки, напальчник и т. д.)"

$ws.Range("A148").Value = "-"
$ws.Range("B148").Value = "мыльно рыльные, бритвенные принадлежности"

$ws.Range("B149").Value = "и пр."

$ws.Range("B152").Value = "11.02.23г."

[void]$ws.Range("A1").Select()
